# Applies odds/value updates to Sheet1 as described in the commit diff.
# Workbook has a single worksheet (Sheet1) with football match odds.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("K2").Value = 1.83
$ws.Range("L2").Value = 4.5
$ws.Range("O2").Value = 1.57
$ws.Range("P2").Value = 2.25
$ws.Range("U2").Value = 2.25
$ws.Range("V2").Value = 1.57
$ws.Range("AK2").Value = 41
$ws.Range("AS2").Value = 351

# Row 4
$ws.Range("G4").Value = 1.44
$ws.Range("H4").Value = 4.33
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 2.25
$ws.Range("Q4").Value = 2.02
$ws.Range("R4").Value = 1.88
$ws.Range("Y4").Value = 9
$ws.Range("AG4").Value = 15
$ws.Range("AH4").Value = 34
$ws.Range("AI4").Value = 21
$ws.Range("AN4").Value = 3.25
$ws.Range("AZ4").Value = 151

# Row 6
$ws.Range("G6").Value = 2.63
$ws.Range("J6").Value = 3.4
$ws.Range("K6").Value = 1.91
$ws.Range("S6").Value = 1.57
$ws.Range("T6").Value = 2.25
$ws.Range("U6").Value = 2.05
$ws.Range("V6").Value = 1.7
$ws.Range("W6").Value = 6.5
$ws.Range("Y6").Value = 11
$ws.Range("AA6").Value = 26
$ws.Range("AC6").Value = 6.5
$ws.Range("AJ6").Value = 29
$ws.Range("AM6").Value = 501
$ws.Range("AS6").Value = 301
$ws.Range("AT6").Value = 2.25
$ws.Range("AY6").Value = 34

# Row 8
$ws.Range("G8").Value = 3.75
$ws.Range("H8").Value = 3.6
$ws.Range("AQ8").Value = 51

# Row 9
$ws.Range("H9").Value = 4.5
$ws.Range("I9").Value = 7.1
$ws.Range("L9").Value = 6.4
$ws.Range("M9").Value = 1.04
$ws.Range("N9").Value = 10
$ws.Range("R9").Value = 2.05
$ws.Range("U9").Value = 1.91
$ws.Range("V9").Value = 1.85
$ws.Range("W9").Value = 6.2
$ws.Range("AC9").Value = 12.5
$ws.Range("AD9").Value = 7.9
$ws.Range("AF9").Value = 70
$ws.Range("AH9").Value = 37
$ws.Range("AI9").Value = 18.5
$ws.Range("AJ9").Value = 120
$ws.Range("AL9").Value = 55
$ws.Range("AW9").Value = 8.5

# Row 10
$ws.Range("Q10").Value = 3.4
$ws.Range("R10").Value = 1.33

# Row 11
$ws.Range("I11").Value = 1.05
$ws.Range("J11").Value = 26
$ws.Range("L11").Value = 1.28
$ws.Range("Q11").Value = 1.28
$ws.Range("U11").Value = 2.9
$ws.Range("V11").Value = 1.36
$ws.Range("W11").Value = 175
$ws.Range("Y11").Value = 250
$ws.Range("AB11").Value = 900
$ws.Range("AE11").Value = 80
$ws.Range("AF11").Value = 450
$ws.Range("AI11").Value = 16.5
$ws.Range("AL11").Value = 65
$ws.Range("AN11").Value = 35
$ws.Range("AO11").Value = 450
$ws.Range("AP11").Value = 200
$ws.Range("AU11").Value = 15.5
$ws.Range("AX11").Value = 3.8
$ws.Range("AY11").Value = 18.5
$ws.Range("AZ11").Value = 6.7
